$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 429, shifting existing rows 429:513 down to 430:514.
$ws.Rows.Item(429).Insert()

# Populate the newly inserted row 429 with the new record's data.
$ws.Cells.Item(429, 1).Value = 5
$ws.Cells.Item(429, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(429, 3).Value = "Maule"
$ws.Cells.Item(429, 4).Value = 45015
$ws.Cells.Item(429, 5).Value = 7
$ws.Cells.Item(429, 6).Value = 100112023
$ws.Cells.Item(429, 7).Value = "Brócoli"
$ws.Cells.Item(429, 8).Value = "Sin especificar"
$ws.Cells.Item(429, 9).Value = "Primera"
$ws.Cells.Item(429, 10).Value = 5000
$ws.Cells.Item(429, 11).Value = 700
$ws.Cells.Item(429, 12).Value = 700
$ws.Cells.Item(429, 13).Value = 700
$ws.Cells.Item(429, 14).Value = "$/unidad"
$ws.Cells.Item(429, 15).Value = "Región del Maule"
$ws.Cells.Item(429, 16).Value = 700
$ws.Cells.Item(429, 17).Value = 1
$ws.Cells.Item(429, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date-time number format used by the
# other records in column D (style index mirrors the existing date cells).
$ws.Cells.Item(429, 4).NumberFormat = $ws.Cells.Item(430, 4).NumberFormat
